$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Number of Promotions" row (row 4). This shifts "Distance from
# Home" up to row 4 and "Number of Dependents" up to row 5.
$ws.Rows.Item(4).Delete()

# Remove the trailing "Number of Dependents" row (now row 5) entirely.
$ws.Rows.Item(5).Delete()
